$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.648.10"
$ws.Range("E2").Value = "  +3.97%  "
$ws.Range("D3").Value = "1.916.95"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.36"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4674"
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.13"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08028"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.014"
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.37"
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.946.18"
$ws.Range("E13").Value = "  +4.68%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.993"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.04"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06586"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.87"
$ws.Range("E20").Value = "  +3.90%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "29.596.31"
$ws.Range("E22").Value = "  +3.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.586"
$ws.Range("E23").Value = "  +4.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.61"
$ws.Range("E24").Value = "  +6.92%  "
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").Value = "2.118.19"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.67"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.92"
$ws.Range("E28").Value = "  +3.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.770"
$ws.Range("E29").Value = "  +8.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.145"
$ws.Range("E30").Value = "  +3.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.58"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.074"
$ws.Range("E32").Value = "  +12.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09470"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.432"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.576"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.408"
$ws.Range("E36").Value = "  +3.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06131"
$ws.Range("E37").Value = "  +1.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02269"
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.423"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.179"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5911"
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.269"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.353"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07524"
$ws.Range("E46").Value = "  +4.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5575"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.19"
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.18"
$ws.Range("E50").Value = "  +2.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2999"
$ws.Range("E51").Value = "  +12.68%  "
